{"js": "// Replace the date line and each \"A\u00d7B=C\" answer cell with its updated value.\n// Every old string below occurs exactly once in the document, so an exact\n// (non-wildcard, case-sensitive) search-and-replace is safe and unambiguous.\nconst replacements = [\n  [\"2025-10-08 Wednesday\", \"2025-10-09 Thursday\"],\n  [\"850\u00d75=4250\", \"431\u00d74=1724\"],\n  [\"349\u00d79=3141\", \"390\u00d77=2730\"],\n  [\"946\u00d75=4730\", \"657\u00d76=3942\"],\n  [\"230\u00d79=2070\", \"747\u00d76=4482\"],\n  [\"811\u00d76=4866\", \"952\u00d72=1904\"],\n  [\"567\u00d76=3402\", \"444\u00d78=3552\"],\n  [\"446\u00d78=3568\", \"460\u00d78=3680\"],\n  [\"126\u00d77=882\", \"480\u00d73=1440\"],\n  [\"564\u00d74=2256\", \"714\u00d79=6426\"],\n  [\"152\u00d74=608\", \"604\u00d77=4228\"],\n  [\"198\u00d77=1386\", \"926\u00d75=4630\"],\n  [\"868\u00d76=5208\", \"341\u00d72=682\"],\n  [\"780\u00d77=5460\", \"948\u00d76=5688\"],\n  [\"267\u00d72=534\", \"421\u00d75=2105\"],\n  [\"817\u00d77=5719\", \"451\u00d77=3157\"],\n  [\"853\u00d79=7677\", \"970\u00d75=4850\"],\n  [\"829\u00d77=5803\", \"555\u00d73=1665\"],\n  [\"789\u00d77=5523\", \"429\u00d79=3861\"],\n  [\"398\u00d72=796\", \"390\u00d77=2730\"],\n  [\"523\u00d73=1569\", \"593\u00d78=4744\"],\n  [\"339\u00d75=1695\", \"866\u00d79=7794\"],\n  [\"611\u00d74=2444\", \"688\u00d76=4128\"],\n  [\"886\u00d76=5316\", \"750\u00d75=3750\"],\n  [\"194\u00d76=1164\", \"645\u00d76=3870\"],\n  [\"453\u00d73=1359\", \"661\u00d77=4627\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the date line and each \"A\u00d7B=C\" answer cell with its updated value.\n# Every old string below occurs exactly once in the document, so a plain\n# Find/Replace (no wildcards) over the whole story is safe and unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"2025-10-08 Wednesday\", \"2025-10-09 Thursday\"),\n  @(\"850\u00d75=4250\", \"431\u00d74=1724\"),\n  @(\"349\u00d79=3141\", \"390\u00d77=2730\"),\n  @(\"946\u00d75=4730\", \"657\u00d76=3942\"),\n  @(\"230\u00d79=2070\", \"747\u00d76=4482\"),\n  @(\"811\u00d76=4866\", \"952\u00d72=1904\"),\n  @(\"567\u00d76=3402\", \"444\u00d78=3552\"),\n  @(\"446\u00d78=3568\", \"460\u00d78=3680\"),\n  @(\"126\u00d77=882\", \"480\u00d73=1440\"),\n  @(\"564\u00d74=2256\", \"714\u00d79=6426\"),\n  @(\"152\u00d74=608\", \"604\u00d77=4228\"),\n  @(\"198\u00d77=1386\", \"926\u00d75=4630\"),\n  @(\"868\u00d76=5208\", \"341\u00d72=682\"),\n  @(\"780\u00d77=5460\", \"948\u00d76=5688\"),\n  @(\"267\u00d72=534\", \"421\u00d75=2105\"),\n  @(\"817\u00d77=5719\", \"451\u00d77=3157\"),\n  @(\"853\u00d79=7677\", \"970\u00d75=4850\"),\n  @(\"829\u00d77=5803\", \"555\u00d73=1665\"),\n  @(\"789\u00d77=5523\", \"429\u00d79=3861\"),\n  @(\"398\u00d72=796\", \"390\u00d77=2730\"),\n  @(\"523\u00d73=1569\", \"593\u00d78=4744\"),\n  @(\"339\u00d75=1695\", \"866\u00d79=7794\"),\n  @(\"611\u00d74=2444\", \"688\u00d76=4128\"),\n  @(\"886\u00d76=5316\", \"750\u00d75=3750\"),\n  @(\"194\u00d76=1164\", \"645\u00d76=3870\"),\n  @(\"453\u00d73=1359\", \"661\u00d77=4627\")\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $rng = $d.Content\n  $rng.Find.ClearFormatting()\n  $rng.Find.Replacement.ClearFormatting()\n  $found = $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n\n  if (-not $found) {\n    throw \"Text not found: $oldText\"\n  }\n}\n"}
